# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets
# to reflect refreshed counts, per the gh-pages data regeneration.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row number -> new value for column F
$updates = @{
    12 = 540
    22 = 83
    23 = 832
    24 = 1360
    25 = 284
    26 = 296
    27 = 186
    33 = 198
    35 = 251
    36 = 1582
    42 = 3330
    44 = 173
    45 = 858
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
